$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2018-12-31 00:00:00"
$ws.Range("O2").Value = 192701973.1
$ws.Range("P2").Value = 76.379704512
$ws.Range("Q2").Value = 3804380968.99
$ws.Range("R2").Value = 1507.9113596402
$ws.Range("S2").Value = 260985242.72
$ws.Range("T2").Value = 103.4445854408
$ws.Range("U2").Value = -24498727.74
$ws.Range("V2").Value = -9.7103602812
$ws.Range("Y2").Value = 173829894.31
$ws.Range("Z2").Value = 68.89953304140001
$ws.Range("AA2").Value = -420581906.02
$ws.Range("AB2").Value = -166.7026091541
$ws.Range("AC2").Value = -252294735.01
$ws.Range("AD2").Value = -236.5737502341
